# Generate Report for Handoff
# Update the UUID-stamped file names, hashes and timestamps across the
# Overview / zh-cn / de-de sheets to reflect the latest handoff run.

$wb = $excel.ActiveWorkbook

$oldGuid = "1af532ad-0c76-4b19-8f41-8894247b8d6f"
$newGuid = "0d1f9c3b-c84d-47f4-9961-88a0548a8b52"

$oldHash = "2f8bf96a04d433e8ca69941a151f74696e34e203"
$newHash = "39b827bd1eac43675dd5082d409c948873d15a3e"

$newMdName = "$newGuid.md"

$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$newHandoffDate = "2016-03-22 08:59:46"
$newZhHandoffDatetime = "2016-03-22 08:59:42"
$newDeHandoffDatetime = "2016-03-22 08:59:46"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffDatetime

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeHandoffDatetime
